$d = $word.ActiveDocument

$pairs = @(
    @("24+30=", "13-1="),
    @("88-21=", "68-47="),
    @("42-1=", "48+10="),
    @("68-34=", "90-15="),
    @("20-16=", "1+56="),
    @("29-5=", "75-13="),
    @("19+30=", "14-0="),
    @("36-7=", "71-21="),
    @("20+40=", "70-62="),
    @("63-36=", "18+17="),
    @("67+0=", "29+17="),
    @("50+5=", "50+12="),
    @("99-49=", "41-13="),
    @("35+38=", "21+11="),
    @("78-32=", "94-30="),
    @("23+46=", "51+47="),
    @("2+23=", "56-23="),
    @("50-13=", "65+3="),
    @("42-38=", "46+39="),
    @("97-40=", "11+12="),
    @("3+82=", "36-17="),
    @("88-16=", "19+59="),
    @("78-22=", "29+33="),
    @("27+64=", "97-33="),
    @("68+14=", "87+0="),
    @("58-26=", "29+34="),
    @("39+13=", "36-33="),
    @("37-35=", "24+43="),
    @("76-76=", "21+62="),
    @("26+2=", "69-43="),
    @("38+61=", "67-34="),
    @("18+38=", "24+73="),
    @("46-32=", "72-30="),
    @("47-38=", "36-3="),
    @("72-59=", "33+60="),
    @("52+12=", "75+9="),
    @("89-47=", "44+12="),
    @("93-6=", "97-36="),
    @("51+3=", "62-41="),
    @("89-71=", "43+53="),
    @("16+32=", "77-26="),
    @("30+54=", "97-20="),
    @("29+38=", "17+28="),
    @("49-35=", "93-37="),
    @("25+61=", "37+18="),
    @("63-45=", "35+12="),
    @("53+12=", "17+13="),
    @("20+49=", "55-24="),
    @("61-8=", "6+14="),
    @("40-10=", "77+16="),
    @("34-11=", "64-33="),
    @("68-25=", "67-51="),
    @("4+25=", "49-8="),
    @("26+49=", "97-26="),
    @("79+4=", "61-6="),
    @("76+2=", "21+26="),
    @("52-38=", "87-7="),
    @("38-8=", "86-61="),
    @("18+45=", "10+57="),
    @("5+30=", "48+44="),
    @("91-82=", "24+37="),
    @("5+35=", "18+65="),
    @("87-83=", "17+59="),
    @("48+17=", "27+27="),
    @("20+37=", "69+28="),
    @("19+33=", "10+83="),
    @("92+5=", "79-4="),
    @("68-28=", "39+10="),
    @("33+16=", "91-7="),
    @("48-7=", "58-42="),
    @("58-40=", "78+0="),
    @("4+41=", "1+85="),
    @("70-12=", "41+2="),
    @("72-32=", "54+35="),
    @("26+24=", "40-26="),
    @("64-30=", "94-83="),
    @("95-82=", "97-59="),
    @("47+24=", "73-73="),
    @("42-36=", "65-12="),
    @("18+11=", "49+21="),
    @("74-52=", "68+11="),
    @("69+29=", "36-4="),
    @("36+35=", "14+66="),
    @("30+13=", "69+18="),
    @("50+7=", "89-24="),
    @("86-81=", "81-40="),
    @("79-22=", "99-42="),
    @("54+32=", "77-59="),
    @("1+49=", "21+73="),
    @("43+51=", "21+72="),
    @("76-24=", "7+33="),
    @("2+78=", "47+4="),
    @("12+9=", "92-27="),
    @("3+47=", "6+24="),
    @("35+62=", "44-11="),
    @("49-13=", "28-5="),
    @("24-11=", "93-76="),
    @("91-66=", "89-5="),
    @("7+44=", "65-39="),
    @("14+46=", "43-15="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Replacements complete: $($pairs.Count) pairs processed"
